{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"212\u00d72=\", \"424\u00d78=\"],\n  [\"470\u00d74=\", \"431\u00d75=\"],\n  [\"695\u00d79=\", \"527\u00d76=\"],\n  [\"358\u00d78=\", \"810\u00d76=\"],\n  [\"134\u00d74=\", \"163\u00d76=\"],\n  [\"318\u00d73=\", \"270\u00d76=\"],\n  [\"576\u00d72=\", \"729\u00d74=\"],\n  [\"867\u00d78=\", \"627\u00d74=\"],\n  [\"376\u00d76=\", \"945\u00d74=\"],\n  [\"858\u00d76=\", \"818\u00d72=\"],\n  [\"906\u00d75=\", \"762\u00d72=\"],\n  [\"709\u00d74=\", \"506\u00d78=\"],\n  [\"752\u00d73=\", \"233\u00d79=\"],\n  [\"974\u00d73=\", \"942\u00d72=\"],\n  [\"179\u00d79=\", \"549\u00d72=\"],\n  [\"157\u00d76=\", \"469\u00d75=\"],\n  [\"272\u00d78=\", \"517\u00d77=\"],\n  [\"735\u00d79=\", \"847\u00d75=\"],\n  [\"832\u00d74=\", \"447\u00d77=\"],\n  [\"964\u00d76=\", \"219\u00d76=\"],\n  [\"707\u00d74=\", \"782\u00d74=\"],\n  [\"879\u00d76=\", \"751\u00d75=\"],\n  [\"380\u00d77=\", \"329\u00d76=\"],\n  [\"712\u00d78=\", \"895\u00d72=\"],\n  [\"415\u00d78=\", \"987\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"212\u00d72=\", \"424\u00d78=\"),\n    @(\"470\u00d74=\", \"431\u00d75=\"),\n    @(\"695\u00d79=\", \"527\u00d76=\"),\n    @(\"358\u00d78=\", \"810\u00d76=\"),\n    @(\"134\u00d74=\", \"163\u00d76=\"),\n    @(\"318\u00d73=\", \"270\u00d76=\"),\n    @(\"576\u00d72=\", \"729\u00d74=\"),\n    @(\"867\u00d78=\", \"627\u00d74=\"),\n    @(\"376\u00d76=\", \"945\u00d74=\"),\n    @(\"858\u00d76=\", \"818\u00d72=\"),\n    @(\"906\u00d75=\", \"762\u00d72=\"),\n    @(\"709\u00d74=\", \"506\u00d78=\"),\n    @(\"752\u00d73=\", \"233\u00d79=\"),\n    @(\"974\u00d73=\", \"942\u00d72=\"),\n    @(\"179\u00d79=\", \"549\u00d72=\"),\n    @(\"157\u00d76=\", \"469\u00d75=\"),\n    @(\"272\u00d78=\", \"517\u00d77=\"),\n    @(\"735\u00d79=\", \"847\u00d75=\"),\n    @(\"832\u00d74=\", \"447\u00d77=\"),\n    @(\"964\u00d76=\", \"219\u00d76=\"),\n    @(\"707\u00d74=\", \"782\u00d74=\"),\n    @(\"879\u00d76=\", \"751\u00d75=\"),\n    @(\"380\u00d77=\", \"329\u00d76=\"),\n    @(\"712\u00d78=\", \"895\u00d72=\"),\n    @(\"415\u00d78=\", \"987\u00d74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}"}
